$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data values for rows 2-6 (data rows 1-5)

# Row 2
$ws.Range("D2").Value = 11067
$ws.Range("E2").Value = 612
$ws.Range("F2").Value = 612
$ws.Range("G2").Value = 791
$ws.Range("H2").Value = 381
$ws.Range("I2").Value = 274
$ws.Range("J2").Value = 107
$ws.Range("K2").Value = 11540
$ws.Range("L2").Value = 2952
$ws.Range("M2").Value = 8588
$ws.Range("N2").Value = 5378
$ws.Range("O2").Value = 3210
$ws.Range("P2").Value = 211
$ws.Range("Q2").Value = 930
$ws.Range("R2").Value = -1250
$ws.Range("S2").Value = 457
$ws.Range("T2").Value = 1226
$ws.Range("U2").Value = -296
$ws.Range("V2").Value = 1442
$ws.Range("W2").Value = 5.53
$ws.Range("X2").Value = 3.44
$ws.Range("Y2").Value = 5.2
$ws.Range("Z2").Value = 3.26
$ws.Range("AA2").Value = 34.38
$ws.Range("AB2").Value = 2409.05
$ws.Range("AC2").Value = 6493
$ws.Range("AD2").Value = 11.5
$ws.Range("AE2").Value = 131613
$ws.Range("AF2").Value = 0.57
$ws.Range("AG2").Value = 2100
$ws.Range("AH2").Value = 2.81
$ws.Range("AI2").Value = 31.28
$ws.Range("AJ2").Value = 4224646

# Row 3
$ws.Range("D3").Value = 10140
$ws.Range("E3").Value = 615
$ws.Range("F3").Value = 615
$ws.Range("G3").Value = 740
$ws.Range("H3").Value = 615
$ws.Range("I3").Value = 340
$ws.Range("J3").Value = 275
$ws.Range("K3").Value = 11712
$ws.Range("L3").Value = 2726
$ws.Range("M3").Value = 8986
$ws.Range("N3").Value = 5570
$ws.Range("O3").Value = 3416
$ws.Range("P3").Value = 211
$ws.Range("Q3").Value = 607
$ws.Range("R3").Value = -760
$ws.Range("S3").Value = -28
$ws.Range("T3").Value = 654
$ws.Range("U3").Value = -46
$ws.Range("V3").Value = 1559
$ws.Range("W3").Value = 6.07
$ws.Range("X3").Value = 6.07
$ws.Range("Y3").Value = 6.21
$ws.Range("Z3").Value = 5.29
$ws.Range("AA3").Value = 30.34
$ws.Range("AB3").Value = 2529.21
$ws.Range("AC3").Value = 8050
$ws.Range("AD3").Value = 7.37
$ws.Range("AE3").Value = 136315
$ws.Range("AF3").Value = 0.44
$ws.Range("AG3").Value = 2100
$ws.Range("AH3").Value = 3.54
$ws.Range("AI3").Value = 25.23
$ws.Range("AJ3").Value = 4224646

# Row 4
$ws.Range("D4").Value = 6965
$ws.Range("E4").Value = 522
$ws.Range("F4").Value = 558
$ws.Range("G4").Value = 633
$ws.Range("H4").Value = 534
$ws.Range("I4").Value = 296
$ws.Range("J4").Value = 239
$ws.Range("K4").Value = 12099
$ws.Range("L4").Value = 2736
$ws.Range("M4").Value = 9362
$ws.Range("N4").Value = 5779
$ws.Range("O4").Value = 3583
$ws.Range("P4").Value = 211
$ws.Range("Q4").Value = 992
$ws.Range("R4").Value = -356
$ws.Range("S4").Value = -206
$ws.Range("T4").Value = 504
$ws.Range("U4").Value = 488
$ws.Range("V4").Value = 1516
$ws.Range("W4").Value = 7.5
$ws.Range("X4").Value = 7.67
$ws.Range("Y4").Value = 5.21
$ws.Range("Z4").Value = 4.49
$ws.Range("AA4").Value = 29.23
$ws.Range("AB4").Value = 2626.73
$ws.Range("AC4").Value = 7000
$ws.Range("AD4").Value = 9.4
$ws.Range("AE4").Value = 141449
$ws.Range("AF4").Value = 0.47
$ws.Range("AG4").Value = 2100
$ws.Range("AH4").Value = 3.19
$ws.Range("AI4").Value = 29.01
$ws.Range("AJ4").Value = 4224646

# Row 5
$ws.Range("D5").Value = 7533
$ws.Range("E5").Value = 364
$ws.Range("F5").Value = 364
$ws.Range("G5").Value = 686
$ws.Range("H5").Value = 584
$ws.Range("I5").Value = 403
$ws.Range("J5").Value = 180
$ws.Range("K5").Value = 15125
$ws.Range("L5").Value = 3220
$ws.Range("M5").Value = 11905
$ws.Range("N5").Value = 6127
$ws.Range("O5").Value = 5777
$ws.Range("P5").Value = 211
$ws.Range("Q5").Value = 622
$ws.Range("R5").Value = -282
$ws.Range("S5").Value = -582
$ws.Range("T5").Value = 361
$ws.Range("U5").Value = 261
$ws.Range("V5").Value = 1013
$ws.Range("W5").Value = 4.84
$ws.Range("X5").Value = 7.75
$ws.Range("Y5").Value = 6.78
$ws.Range("Z5").Value = 4.29
$ws.Range("AA5").Value = 27.05
$ws.Range("AB5").Value = 2775.88
$ws.Range("AC5").Value = 9548
$ws.Range("AD5").Value = 6.86
$ws.Range("AE5").Value = 149964
$ws.Range("AF5").Value = 0.44
$ws.Range("AG5").Value = 3850
$ws.Range("AH5").Value = 5.88
$ws.Range("AI5").Value = 39
$ws.Range("AJ5").Value = 4224646

# Row 6
$ws.Range("D6").Value = 10224
$ws.Range("E6").Value = 338
$ws.Range("F6").Value = 338
$ws.Range("G6").Value = 269
$ws.Range("H6").Value = 147
$ws.Range("I6").Value = 37
$ws.Range("K6").Value = 14166
$ws.Range("L6").Value = 2907
$ws.Range("M6").Value = 11259
$ws.Range("N6").Value = 6058
$ws.Range("P6").Value = 211
$ws.Range("Q6").Value = 434
$ws.Range("R6").Value = -1164
$ws.Range("S6").Value = 343
$ws.Range("T6").Value = 609
$ws.Range("U6").Value = -175
$ws.Range("V6").Value = 1256
$ws.Range("W6").Value = 3.31
$ws.Range("X6").Value = 1.44
$ws.Range("Y6").Value = 0.6
$ws.Range("Z6").Value = 1.01
$ws.Range("AA6").Value = 25.82
$ws.Range("AB6").Value = 2781.04
$ws.Range("AC6").Value = 871
$ws.Range("AD6").Value = 70.23
$ws.Range("AE6").Value = 148271
$ws.Range("AF6").Value = 0.41
$ws.Range("AG6").Value = 2600
$ws.Range("AH6").Value = 4.25
$ws.Range("AI6").Value = 288.57
$ws.Range("AJ6").Value = 4224646

# Clear out data cells D:AJ for rows 7, 8, 9 (keep only A, B, C)
$ws.Range("D7:AJ9").ClearContents()
